$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.307.89"
$ws.Cells.Item(2, 5).Value = "  +1.25%  "

$ws.Cells.Item(3, 4).Value = "1.809.96"
$ws.Cells.Item(3, 5).Value = "  +3.41%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "337.98"
$ws.Cells.Item(5, 5).Value = "  +0.65%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9993"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4658"
$ws.Cells.Item(7, 5).Value = "  +20.96%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3809"
$ws.Cells.Item(8, 5).Value = "  +11.88%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "45.43"
$ws.Cells.Item(9, 5).Value = "  -0.73%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.160"
$ws.Cells.Item(10, 5).Value = "  +3.74%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07660"
$ws.Cells.Item(11, 5).Value = "  +6.06%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "22.56"
$ws.Cells.Item(12, 5).Value = "  -0.41%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.001"
$ws.Cells.Item(13, 5).Value = "  -0.17%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.352"
$ws.Cells.Item(14, 5).Value = "  +3.00%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.463"
$ws.Cells.Item(15, 5).Value = "  +4.87%  "

$ws.Cells.Item(16, 4).Value = "1.808.01"
$ws.Cells.Item(16, 5).Value = "  +3.18%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001097"
$ws.Cells.Item(17, 5).Value = "  +3.41%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.06718"
$ws.Cells.Item(18, 5).Value = "  +1.60%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "82.00"
$ws.Cells.Item(19, 5).Value = "  +3.58%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.9992"
$ws.Cells.Item(20, 5).Value = "  -0.05%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.51"
$ws.Cells.Item(21, 5).Value = "  +4.66%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.442"
$ws.Cells.Item(22, 5).Value = "  +4.15%  "

$ws.Cells.Item(23, 4).Value = "28.297.24"
$ws.Cells.Item(23, 5).Value = "  +1.14%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.92"
$ws.Cells.Item(24, 5).Value = "  +2.26%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.412"
$ws.Cells.Item(25, 5).Value = "  +0.62%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "20.85"
$ws.Cells.Item(26, 5).Value = "  +4.99%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "154.43"
$ws.Cells.Item(27, 5).Value = "  +0.61%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.382"
$ws.Cells.Item(28, 5).Value = "  +3.67%  "

$ws.Cells.Item(29, 4).Value = "2.014.00"
$ws.Cells.Item(29, 5).Value = "  +3.25%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "133.64"
$ws.Cells.Item(30, 5).Value = "  +1.78%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.262"
$ws.Cells.Item(31, 5).Value = "  +0.07%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.033"
$ws.Cells.Item(32, 5).Value = "  +0.19%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.09604"
$ws.Cells.Item(33, 5).Value = "  +8.86%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.880"
$ws.Cells.Item(34, 5).Value = "  +0.59%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.2277"
$ws.Cells.Item(35, 5).Value = "  +8.66%  "

$ws.Cells.Item(36, 2).Value = "Aptos"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "12.16"
$ws.Cells.Item(36, 5).Value = "  -0.27%  "

$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.06390"
$ws.Cells.Item(37, 5).Value = "  +4.18%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02359"
$ws.Cells.Item(38, 5).Value = "  +3.32%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "5.288"
$ws.Cells.Item(39, 5).Value = "  +3.01%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.6665"
$ws.Cells.Item(40, 5).Value = "  +1.71%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.241"
$ws.Cells.Item(41, 5).Value = "  +2.91%  "

$ws.Cells.Item(42, 5).Value = "  -2.74%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "8.340"
$ws.Cells.Item(43, 5).Value = "  +4.09%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "14.15"
$ws.Cells.Item(44, 5).Value = "  +3.30%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.9993"
$ws.Cells.Item(45, 5).Value = "  +0.01%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.6162"
$ws.Cells.Item(46, 5).Value = "  +2.14%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.861"
$ws.Cells.Item(47, 5).Value = "  +0.58%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "131.42"
$ws.Cells.Item(48, 5).Value = "  +3.80%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.046"
$ws.Cells.Item(49, 5).Value = "  +2.12%  "

$ws.Cells.Item(50, 2).Value = "EOS"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.182"
$ws.Cells.Item(50, 5).Value = "  +0.92%  "

$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.07162"
$ws.Cells.Item(51, 5).Value = "  +2.57%  "
